# Update column G ("K") values in the active worksheet to reflect the
# regenerated strike-count (K) data, replacing the old Strike# derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 5
    20 = 5
    21 = 0
    22 = 0
    23 = 1
    24 = 1
    25 = 1
    26 = 2
    27 = 2
    28 = 1
    29 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
